$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Copy()
$ws.Range("A67").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A67").Value = 46016
$ws.Range("B67").Value = 144
$ws.Range("C67").Value = 157
$ws.Range("D67").Value = 147
